$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear stale cells that must become empty
$ws.Cells.Item(2,13).ClearContents()
$ws.Cells.Item(9,29).ClearContents()
$ws.Cells.Item(11,13).ClearContents()
$ws.Cells.Item(11,36).ClearContents()
$ws.Cells.Item(11,37).ClearContents()
$ws.Cells.Item(11,41).ClearContents()
$ws.Cells.Item(13,13).ClearContents()
$ws.Cells.Item(13,36).ClearContents()
$ws.Cells.Item(13,37).ClearContents()
$ws.Cells.Item(13,41).ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = 111814135
$ws.Cells.Item(2,2).Value = 90087
$ws.Cells.Item(2,4).Value = "LC"
$ws.Cells.Item(2,5).Value = 3298
$ws.Cells.Item(2,6).Value = "Trådticka"
$ws.Cells.Item(2,7).Value = "Climacocystis borealis"
$ws.Cells.Item(2,8).Value = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(2,17).Value = 540661.0419420782
$ws.Cells.Item(2,18).Value = 7247564.172119373
# Row 3
$ws.Cells.Item(3,1).Value = 111814152
$ws.Cells.Item(3,17).Value = 540661.0419420782
$ws.Cells.Item(3,18).Value = 7247564.172119373
# Row 4
$ws.Cells.Item(4,1).Value = 111813872
$ws.Cells.Item(4,2).Value = 56398
$ws.Cells.Item(4,4).Value = "NT"
$ws.Cells.Item(4,5).Value = 100109
$ws.Cells.Item(4,6).Value = "Tretåig hackspett"
$ws.Cells.Item(4,7).Value = "Picoides tridactylus"
$ws.Cells.Item(4,8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(4,13).Value = "färska spår"
$ws.Cells.Item(4,17).Value = 540557.5018987871
$ws.Cells.Item(4,18).Value = 7247552.715308581
# Row 5
$ws.Cells.Item(5,1).Value = 111814303
$ws.Cells.Item(5,2).Value = 90087
$ws.Cells.Item(5,4).Value = "LC"
$ws.Cells.Item(5,5).Value = 3298
$ws.Cells.Item(5,6).Value = "Trådticka"
$ws.Cells.Item(5,7).Value = "Climacocystis borealis"
$ws.Cells.Item(5,8).Value = "(Fr.) Kotl. & Pouzar"
$ws.Cells.Item(5,17).Value = 540600.641023421
$ws.Cells.Item(5,18).Value = 7247517.393825463
# Row 6
$ws.Cells.Item(6,1).Value = 111813938
$ws.Cells.Item(6,2).Value = 89423
$ws.Cells.Item(6,5).Value = 5432
$ws.Cells.Item(6,6).Value = "Granticka"
$ws.Cells.Item(6,7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(6,8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(6,17).Value = 540654.849203686
$ws.Cells.Item(6,18).Value = 7247498.096959669
# Row 7
$ws.Cells.Item(7,1).Value = 111813785
$ws.Cells.Item(7,2).Value = 89405
$ws.Cells.Item(7,5).Value = 1202
$ws.Cells.Item(7,6).Value = "Ullticka"
$ws.Cells.Item(7,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(7,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(7,17).Value = 540570.9514120822
$ws.Cells.Item(7,18).Value = 7247577.960198429
# Row 8
$ws.Cells.Item(8,1).Value = 111813707
$ws.Cells.Item(8,2).Value = 56398
$ws.Cells.Item(8,5).Value = 100109
$ws.Cells.Item(8,6).Value = "Tretåig hackspett"
$ws.Cells.Item(8,7).Value = "Picoides tridactylus"
$ws.Cells.Item(8,8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(8,13).Value = "färska spår"
$ws.Cells.Item(8,17).Value = 540647.037727406
$ws.Cells.Item(8,18).Value = 7247579.013394679
$ws.Cells.Item(8,36).Value = "gran"
$ws.Cells.Item(8,37).Value = "Picea abies"
$ws.Cells.Item(8,41).Value = "Picea abies"
# Row 9
$ws.Cells.Item(9,1).Value = 111814212
$ws.Cells.Item(9,2).Value = 89405
$ws.Cells.Item(9,4).Value = "NT"
$ws.Cells.Item(9,5).Value = 1202
$ws.Cells.Item(9,6).Value = "Ullticka"
$ws.Cells.Item(9,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(9,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(9,17).Value = 540635.9369002836
$ws.Cells.Item(9,18).Value = 7247595.565451766
# Row 10
$ws.Cells.Item(10,1).Value = 111813745
$ws.Cells.Item(10,2).Value = 56398
$ws.Cells.Item(10,5).Value = 100109
$ws.Cells.Item(10,6).Value = "Tretåig hackspett"
$ws.Cells.Item(10,7).Value = "Picoides tridactylus"
$ws.Cells.Item(10,8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(10,13).Value = "färska spår"
$ws.Cells.Item(10,17).Value = 540568.950047517
$ws.Cells.Item(10,18).Value = 7247601.73830481
$ws.Cells.Item(10,36).Value = "gran"
$ws.Cells.Item(10,37).Value = "Picea abies"
$ws.Cells.Item(10,41).Value = "Picea abies"
# Row 11
$ws.Cells.Item(11,1).Value = 111814119
$ws.Cells.Item(11,2).Value = 89423
$ws.Cells.Item(11,5).Value = 5432
$ws.Cells.Item(11,6).Value = "Granticka"
$ws.Cells.Item(11,7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(11,8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(11,17).Value = 540683.0369185829
$ws.Cells.Item(11,18).Value = 7247576.171207689
# Row 12
$ws.Cells.Item(12,1).Value = 111814047
$ws.Cells.Item(12,17).Value = 540633.6855369165
$ws.Cells.Item(12,18).Value = 7247516.598344535
$ws.Cells.Item(12,29).Value = "rikligt"
# Row 13
$ws.Cells.Item(13,1).Value = 111813975
$ws.Cells.Item(13,2).Value = 89423
$ws.Cells.Item(13,5).Value = 5432
$ws.Cells.Item(13,6).Value = "Granticka"
$ws.Cells.Item(13,7).Value = "Porodaedalea chrysoloma"
$ws.Cells.Item(13,8).Value = "(Fr.) Fiasson & Niemelä"
$ws.Cells.Item(13,17).Value = 540643.7191088985
$ws.Cells.Item(13,18).Value = 7247516.737328541
# Row 14
$ws.Cells.Item(14,1).Value = 111825158
$ws.Cells.Item(14,2).Value = 89405
$ws.Cells.Item(14,4).Value = "NT"
$ws.Cells.Item(14,5).Value = 1202
$ws.Cells.Item(14,6).Value = "Ullticka"
$ws.Cells.Item(14,7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(14,8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(14,29).Value = "Med ulltickeporing"
# Row 15
$ws.Cells.Item(15,1).Value = 111825245
$ws.Cells.Item(15,2).Value = 89745
$ws.Cells.Item(15,4).Value = "VU"
$ws.Cells.Item(15,5).Value = 2062
$ws.Cells.Item(15,6).Value = "Ulltickeporing"
$ws.Cells.Item(15,7).Value = "Skeletocutis brevispora"
$ws.Cells.Item(15,8).Value = "Niemelä"
$ws.Cells.Item(15,29).Value = "Färskt exemplar. Kollekt tog och torkades, gulnade."
$ws.Cells.Item(15,36).Value = "ullticka"
$ws.Cells.Item(15,37).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(15,41).Value = "Phellinidium ferrugineofuscum"
# Row 17
$ws.Cells.Item(17,1).Value = 111825098
$ws.Cells.Item(17,2).Value = 89845
$ws.Cells.Item(17,5).Value = 1209
$ws.Cells.Item(17,6).Value = "Rynkskinn"
$ws.Cells.Item(17,7).Value = "Phlebia centrifuga"
$ws.Cells.Item(17,8).Value = "P.Karst."
$ws.Cells.Item(17,29).Value = "Färskt. På granlåga med minst 25 rosentickor, ullticka och ulltickeporing.  Någon gul slemsvamp? på rynkskinnet"
$ws.Cells.Item(17,36).Value = "gran"
$ws.Cells.Item(17,37).Value = "Picea abies"
$ws.Cells.Item(17,41).Value = "Picea abies"
# Row 18
$ws.Cells.Item(18,1).Value = 111905851
$ws.Cells.Item(18,2).Value = 56398
$ws.Cells.Item(18,3).Value = "Ovaliderad"
$ws.Cells.Item(18,4).Value = "NT"
$ws.Cells.Item(18,5).Value = 100109
$ws.Cells.Item(18,6).Value = "Tretåig hackspett"
$ws.Cells.Item(18,7).Value = "Picoides tridactylus"
$ws.Cells.Item(18,8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(18,13).Value = "färska spår"
$ws.Cells.Item(18,16).Value = "Matsdal, granskog, Ås lm"
$ws.Cells.Item(18,17).Value = 540642.8745156997
$ws.Cells.Item(18,18).Value = 7247577.702774134
$ws.Cells.Item(18,19).Value = 10
$ws.Cells.Item(18,20).Value = "Västerbotten"
$ws.Cells.Item(18,21).Value = "Vilhelmina"
$ws.Cells.Item(18,22).Value = "Åsele lappmark"
$ws.Cells.Item(18,23).Value = "Vilhelmina"
$ws.Cells.Item(18,25).NumberFormat = "@"
$ws.Cells.Item(18,25).Value = "2023-08-13"
$ws.Cells.Item(18,26).NumberFormat = "@"
$ws.Cells.Item(18,26).Value = "00:00"
$ws.Cells.Item(18,27).NumberFormat = "@"
$ws.Cells.Item(18,27).Value = "2023-08-13"
$ws.Cells.Item(18,28).NumberFormat = "@"
$ws.Cells.Item(18,28).Value = "00:00"
$ws.Cells.Item(18,30).Value = $false
$ws.Cells.Item(18,31).Value = $false
$ws.Cells.Item(18,33).Value = $false
$ws.Cells.Item(18,49).Value = "Roger Olofsson"
$ws.Cells.Item(18,50).Value = "Roger Olofsson"
